# Atualização de bases das ligas, do dia: 17-05-2024 às 13:59
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the data (columns B..AB) between paired rows whose fixture
# records were reordered in the source feed. Column A (the running id)
# stays put; everything else (match id, teams, scores, odds) swaps. ---

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AB$rowA")
    $rangeB = $ws.Range("B$rowB`:AB$rowB")
    $valsA = $rangeA.Value()
    $valsB = $rangeB.Value()
    $rangeA.Value = $valsB
    $rangeB.Value = $valsA
}

Swap-Rows 130 131
Swap-Rows 142 144
Swap-Rows 143 145

# --- Odds refresh for upcoming fixtures (rows 242, 243, 244, 249) ---

$ws.Range("M242").Value = 1.615
$ws.Range("N242").Value = 3.6
$ws.Range("O242").Value = 5.5
$ws.Range("Q242").Value = 1.8
$ws.Range("R242").Value = 2

$ws.Range("M243").Value = 3.2
$ws.Range("O243").Value = 2.25
$ws.Range("Q243").Value = 1.85
$ws.Range("R243").Value = 1.95
$ws.Range("T243").Value = 1.95
$ws.Range("U243").Value = 1.85

$ws.Range("M244").Value = 1.75
$ws.Range("N244").Value = 3.5
$ws.Range("O244").Value = 4.333
$ws.Range("P244").Value = -0.75
$ws.Range("Q244").Value = 2
$ws.Range("R244").Value = 1.8
$ws.Range("T244").Value = 1.975
$ws.Range("U244").Value = 1.825

$ws.Range("M249").Value = 2.1
$ws.Range("N249").Value = 3.1
$ws.Range("O249").Value = 3.6
$ws.Range("P249").Value = -0.25
$ws.Range("Q249").Value = 1.775
$ws.Range("R249").Value = 2.025
